$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("T17").Value = "ID"
$ws.Range("U17").Value = "Slope"
$ws.Range("V17").Value = "δ"

$ws.Range("T18").Value = 307
$ws.Range("U18").Formula = "-0.2274"
$ws.Range("V18").Formula = "=-U18*LN(10)"
$ws.Range("T19").Value = 307
$ws.Range("U19").Formula = "-0.5197"
$ws.Range("V19").Formula = "=-U19*LN(10)"
$ws.Range("U18:U19,V18:V19").NumberFormat = "0.000"

$ws.Range("U25").Value = "ID"
$ws.Range("V25").Value = "Day"
$ws.Range("W25").Value = "V"

$ws.Range("U26").Value = 112
$ws.Range("V26").Value = 0
$ws.Range("W26").Value = 0

$ws.Range("U27").Value = 112
$ws.Range("V27").Value = 1
$ws.Range("W27").Value = 6667
$ws.Range("U28").Value = 112
$ws.Range("V28").Value = 2
$ws.Range("W28").Value = 171203
$ws.Range("W27:W28").Interior.Color = 65535

$ws.Range("U29").Value = 112
$ws.Range("V29").Value = 3
$ws.Range("W29").Value = 786338
$ws.Range("U30").Value = 112
$ws.Range("V30").Value = 4
$ws.Range("W30").Value = 1397797
$ws.Range("U31").Value = 112
$ws.Range("V31").Value = 5
$ws.Range("W31").Value = 284
$ws.Range("U32").Value = 112
$ws.Range("V32").Value = 6
$ws.Range("W32").Value = 0
$ws.Range("U33").Value = 112
$ws.Range("V33").Value = 7
$ws.Range("W33").Value = 0

Write-Host "done"
